$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to hold $value as text, even if it looks numeric,
    # without permanently changing the cell's style/number format.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = $origStyle
}

$updates = @(
    @{ Row = 2;  D = "64.416.35";  E = "  +4.34%  " },
    @{ Row = 3;  D = "2.964.57";   E = "  +2.33%  " },
    @{ Row = 4;  D = $null;        E = "  -0.28%  " },
    @{ Row = 5;  D = "579.20";     E = "  +0.72%  " },
    @{ Row = 6;  D = $null;        E = "  +5.45%  " },
    @{ Row = 7;  D = $null;        E = "  +0.03%  " },
    @{ Row = 8;  D = "2.959.17";   E = "  +2.17%  " },
    @{ Row = 9;  D = $null;        E = "  +0.99%  " },
    @{ Row = 10; D = $null;        E = "  +4.57%  " },
    @{ Row = 11; D = $null;        E = "  +1.43%  " },
    @{ Row = 12; D = "0.444";      E = "  +2.84%  " },
    @{ Row = 13; D = $null;        E = "  +2.73%  " },
    @{ Row = 14; D = "34.24";      E = "  +5.16%  " },
    @{ Row = 15; D = $null;        E = "  +0.79%  " },
    @{ Row = 16; D = "3.460.73";   E = "  +2.44%  " },
    @{ Row = 17; D = "64.338.43";  E = "  +4.06%  " },
    @{ Row = 18; D = "6.91";       E = "  +4.16%  " },
    @{ Row = 19; D = "2.963.03";   E = "  +0.92%  " },
    @{ Row = 20; D = "448.08";     E = "  +3.13%  " },
    @{ Row = 21; D = "13.55";      E = "  +1.97%  " },
    @{ Row = 22; D = "0.672";      E = "  +2.27%  " },
    @{ Row = 23; D = "7.14";       E = "  +3.03%  " },
    @{ Row = 24; D = "80.51";      E = "  +1.03%  " },
    @{ Row = 25; D = "10.94";      E = "  +7.93%  " },
    @{ Row = 26; D = "12.24";      E = "  +3.15%  " },
    @{ Row = 27; D = $null;        E = "  +6.95%  " },
    @{ Row = 28; D = $null;        E = "  -0.02%  " },
    @{ Row = 29; D = "7.60";       E = "  +8.25%  " },
    @{ Row = 30; D = "0.0000109";  E = "  +0.76%  " },
    @{ Row = 31; D = "2.15";       E = "  +2.77%  " },
    @{ Row = 32; D = "2.55";       E = "  +1.02%  " },
    @{ Row = 33; D = $null;        E = "  +2.78%  " },
    @{ Row = 34; D = "26.48";      E = "  +3.08%  " },
    @{ Row = 35; D = "0.999";      E = "  -0.44%  " },
    @{ Row = 36; D = "0.969";      E = "  +0.67%  " },
    @{ Row = 37; D = "2.13";       E = "  +8.79%  " },
    @{ Row = 38; D = "5.60";       E = "  +2.65%  " },
    @{ Row = 39; D = "3.02";       E = "  -0.57%  " },
    @{ Row = 40; D = "49.03";      E = "  -0.12%  " },
    @{ Row = 41; D = "43.94";      E = "  +14.39%  " },
    @{ Row = 42; D = $null;        E = "  +2.62%  " },
    @{ Row = 43; D = "0.295";      E = "  +9.66%  " },
    @{ Row = 44; D = "8.30";       E = "  +0.23%  " },
    @{ Row = 45; D = "379.42";     E = "  +10.99%  " },
    @{ Row = 46; D = "0.0350";     E = "  +4.65%  " },
    @{ Row = 47; D = "2.749.93";   E = "  +2.54%  " },
    @{ Row = 48; D = "134.25";     E = "  -0.06%  " },
    @{ Row = 49; D = $null;        E = "  +0.01%  " },
    @{ Row = 50; D = $null;        E = "  +2.26%  " },
    @{ Row = 51; D = $null;        E = "  +8.52%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    }
    Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
}
